# Add a new "EventID" column at the front of the Events sheet (sheet1),
# populate it with sequential IDs, and make "Events" the active/selected
# sheet with cell C7 selected (Employee sheet selection stays at C10 but
# is no longer the active tab).

$wb = $excel.ActiveWorkbook
$wsEvents = $wb.Worksheets.Item("Events")
$wsEmployee = $wb.Worksheets.Item("Employee")

# Insert a brand-new column before column A; this shifts the existing
# data (and column width formatting) from A:F to B:G automatically.
$wsEvents.Columns.Item(1).EntireColumn.Insert()

# Header cell for the newly inserted column, styled like the other
# header cells (bold header style).
$wsEvents.Range("A1").Value = "EventID"
$wsEvents.Range("A1").Font.Bold = $true

# Sequential EventID values for the six data rows.
$wsEvents.Range("A2").Value = 1
$wsEvents.Range("A3").Value = 2
$wsEvents.Range("A4").Value = 3
$wsEvents.Range("A5").Value = 4
$wsEvents.Range("A6").Value = 5
$wsEvents.Range("A7").Value = 6

# Make "Events" the active sheet/tab, with C7 selected.
$wsEvents.Activate()
$wsEvents.Range("C7").Select() | Out-Null

# Employee sheet keeps its own last-known selection but is no longer
# the active tab.
$wsEmployee.Range("C10").Select() | Out-Null
$wsEvents.Activate()
